# Update bus voltage magnitude results for the 380 kV case (rows 2-25)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.03508387794875
$ws.Range("D2").Value = 1.038959112610628
$ws.Range("E2").Value = 0.992614727750844
$ws.Range("F2").Value = 1.044443363667182
$ws.Range("I2").Value = 1.038340563394524
$ws.Range("J2").Value = 1.040199588978427
$ws.Range("K2").Value = 1.041745662628496
$ws.Range("L2").Value = 0.9955398523335997
$ws.Range("M2").Value = 1.047214398958442
$ws.Range("N2").Value = 1.017324353309392

# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.036061397472625
$ws.Range("D3").Value = 1.039708767190242
$ws.Range("E3").Value = 0.9936372048519299
$ws.Range("F3").Value = 1.045656601634446
$ws.Range("I3").Value = 1.038608154690166
$ws.Range("J3").Value = 1.040820367756812
$ws.Range("K3").Value = 1.042305597133657
$ws.Range("L3").Value = 0.9963617723202687
$ws.Range("M3").Value = 1.048237827325504
$ws.Range("N3").Value = 1.017533305539295

# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.036693845215637
$ws.Range("D4").Value = 1.040193756117787
$ws.Range("E4").Value = 0.9942998659930998
$ws.Range("F4").Value = 1.046441998090052
$ws.Range("I4").Value = 1.038780110454255
$ws.Range("J4").Value = 1.041221381246135
$ws.Range("K4").Value = 1.042667161295123
$ws.Range("L4").Value = 0.9968940712668347
$ws.Range("M4").Value = 1.048899839723459
$ws.Range("N4").Value = 1.017668209137904

# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.036959708800484
$ws.Range("D5").Value = 1.040397623850994
$ws.Range("E5").Value = 0.994578699834602
$ws.Range("F5").Value = 1.046772263392677
$ws.Range("I5").Value = 1.038852114740509
$ws.Range("J5").Value = 1.041389806248352
$ws.Range("K5").Value = 1.042818982786611
$ws.Range("L5").Value = 0.9971179600053012
$ws.Range("M5").Value = 1.049178099010734
$ws.Range("N5").Value = 1.017724850055826

# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.037004347442381
$ws.Range("D6").Value = 1.040431852855267
$ws.Range("E6").Value = 0.994625531979634
$ws.Range("F6").Value = 1.046827721375165
$ws.Range("I6").Value = 1.038864187806549
$ws.Range("J6").Value = 1.041418076095889
$ws.Range("K6").Value = 1.042844463713856
$ws.Range("L6").Value = 0.9971555583673455
$ws.Range("M6").Value = 1.049224817011082
$ws.Range("N6").Value = 1.017734356056329

# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.036697397766449
$ws.Range("D7").Value = 1.040196480294107
$ws.Range("E7").Value = 0.994303590798249
$ws.Range("F7").Value = 1.046446410777685
$ws.Range("I7").Value = 1.038781073702686
$ws.Range("J7").Value = 1.041223632382416
$ws.Range("K7").Value = 1.042669190649315
$ws.Range("L7").Value = 0.9968970624462089
$ws.Range("M7").Value = 1.048903558036289
$ws.Range("N7").Value = 1.017668966261893

# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.035414250295467
$ws.Range("D8").Value = 1.039212479679644
$ws.Range("E8").Value = 0.9929600610674297
$ws.Range("F8").Value = 1.044853311027014
$ws.Range("I8").Value = 1.038431244202496
$ws.Range("J8").Value = 1.04040952306259
$ws.Range("K8").Value = 1.041935050407854
$ws.Range("L8").Value = 0.9958175282591056
$ws.Range("M8").Value = 1.04756031604846
$ws.Range("N8").Value = 1.017395032235279

# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.033152624045191
$ws.Range("D9").Value = 1.037477899758061
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.04204871499251
$ws.Range("I9").Value = 1.037805663767998
$ws.Range("J9").Value = 1.038969821716493
$ws.Range("K9").Value = 1.040635663961925
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.045191686520374
$ws.Range("N9").Value = 1.016910015151339

# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.03164449700387
$ws.Range("D10").Value = 1.03632111344021
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.040180715071189
$ws.Range("I10").Value = 1.037382474638244
$ws.Range("J10").Value = 1.038006577479572
$ws.Range("K10").Value = 1.039765565106321
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.043611442644321
$ws.Range("N10").Value = 1.016585124765853

# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.030991370670861
$ws.Range("D11").Value = 1.035820123949057
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.039372248892747
$ws.Range("I11").Value = 1.037197774928161
$ws.Range("J11").Value = 1.037588666866625
$ws.Range("K11").Value = 1.039387894462973
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.042926895815578
$ws.Range("N11").Value = 1.016444078423518

# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.030748755602423
$ws.Range("D12").Value = 1.035634020570872
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.039072006018485
$ws.Range("I12").Value = 1.037128950510561
$ws.Range("J12").Value = 1.03743331316808
$ws.Range("K12").Value = 1.039247473888574
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.042672580052871
$ws.Range("N12").Value = 1.016391632527559

# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.030800798034816
$ws.Range("D13").Value = 1.035673940984244
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.039136406552399
$ws.Range("I13").Value = 1.037143723478893
$ws.Range("J13").Value = 1.037466642633339
$ws.Range("K13").Value = 1.03927760077023
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.042727133679565
$ws.Range("N13").Value = 1.016402884838022

# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.030971316338298
$ws.Range("D14").Value = 1.035804740860057
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.039347429557542
$ws.Range("I14").Value = 1.03719209034201
$ws.Range("J14").Value = 1.037575827792247
$ws.Range("K14").Value = 1.039376290045965
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.042905874910148
$ws.Range("N14").Value = 1.016439744350995

# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.031076376315618
$ws.Range("D15").Value = 1.03588532919528
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.039477455375508
$ws.Range("I15").Value = 1.037221861780863
$ws.Range("J15").Value = 1.037643084062605
$ws.Range("K15").Value = 1.039437077632242
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.04301599734584
$ws.Range("N15").Value = 1.016462447427299

# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.031687840775652
$ws.Range("D16").Value = 1.036354360534222
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.040234378429987
$ws.Range("I16").Value = 1.037394701875465
$ws.Range("J16").Value = 1.038034295567337
$ws.Range("K16").Value = 1.039790610659134
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.043656867559004
$ws.Range("N16").Value = 1.016594477828988

# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.032071370038945
$ws.Range("D17").Value = 1.036648546822776
$ws.Range("E17").Value = 0.989476357848556
$ws.Range("F17").Value = 1.040709279897933
$ws.Range("I17").Value = 1.037502730102205
$ws.Range("J17").Value = 1.038279472827278
$ws.Range("K17").Value = 1.04001212849921
$ws.Range("L17").Value = 0.9930127773699352
$ws.Range("M17").Value = 1.044058790019722
$ws.Range("N17").Value = 1.016677198904588

# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.03229506679047
$ws.Range("D18").Value = 1.036820131667337
$ws.Range("E18").Value = 0.9897087662937556
$ws.Range("F18").Value = 1.040986319606872
$ws.Range("I18").Value = 1.037565600651924
$ws.Range("J18").Value = 1.038422401502836
$ws.Range("K18").Value = 1.040141248093382
$ws.Range("L18").Value = 0.9932001317071769
$ws.Range("M18").Value = 1.044293196593248
$ws.Range("N18").Value = 1.016725413325474

# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.032371340018328
$ws.Range("D19").Value = 1.036878636117894
$ws.Range("E19").Value = 0.9897880325774034
$ws.Range("F19").Value = 1.041080789420761
$ws.Range("I19").Value = 1.037587014054688
$ws.Range("J19").Value = 1.038471123074618
$ws.Range("K19").Value = 1.040185259604912
$ws.Range("L19").Value = 0.9932640239640975
$ws.Range("M19").Value = 1.044373118455728
$ws.Range("N19").Value = 1.016741847192933

# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.032030221950385
$ws.Range("D20").Value = 1.036616984349882
$ws.Range("E20").Value = 0.9894336180360679
$ws.Range("F20").Value = 1.040658323571659
$ws.Range("I20").Value = 1.037491154226839
$ws.Range("J20").Value = 1.038253175806103
$ws.Range("K20").Value = 1.039988370833947
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.044015670445726
$ws.Range("N20").Value = 1.016668327375245

# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.030921103360322
$ws.Range("D21").Value = 1.035766223927819
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.039285286965141
$ws.Range("I21").Value = 1.037177853543878
$ws.Range("J21").Value = 1.037543678881939
$ws.Range("K21").Value = 1.039347232291427
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.042853241310912
$ws.Range("N21").Value = 1.016428891660855

# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.030223670343129
$ws.Range("D22").Value = 1.035231238768459
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.038422336241275
$ws.Range("I22").Value = 1.036979603529615
$ws.Range("J22").Value = 1.037096877863516
$ws.Range("K22").Value = 1.038943331161687
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.042122117542032
$ws.Range("N22").Value = 1.016278030726163

# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.030593401026028
$ws.Range("D23").Value = 1.035514851795982
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.038879771671883
$ws.Range("I23").Value = 1.037084819552207
$ws.Range("J23").Value = 1.037333803029725
$ws.Range("K23").Value = 1.039157521803859
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.042509724882559
$ws.Range("N23").Value = 1.016358035096226

# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.032048815029212
$ws.Range("D24").Value = 1.036631246100695
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.040681348425204
$ws.Range("I24").Value = 1.037496385300666
$ws.Range("J24").Value = 1.038265058542253
$ws.Range("K24").Value = 1.039999106172122
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.04403515441135
$ws.Range("N24").Value = 1.016672336146791

# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.03373737479797
$ws.Range("D25").Value = 1.037926403175245
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.042773461542064
$ws.Range("I25").Value = 1.037968473253529
$ws.Range("J25").Value = 1.039342626684479
$ws.Range("K25").Value = 1.040972264739
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.045804235290678
$ws.Range("N25").Value = 1.017035676684767
